$wb = $excel.ActiveWorkbook

# --- Sheet1: rename to login1, set login data ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "login1"

$ws1.Cells.Clear()

$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "password"
$ws1.Range("B2").Value = "pass123"
$ws1.Range("A2").Value = "id_123"

# --- Sheet2: add new sheet named GTM (placed after login1) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "GTM"

$ws2.Range("A1").Value = "first name "
$ws2.Range("B1").Value = "lastname"
$ws2.Range("C1").Value = "email"
$ws2.Range("D1").Value = "phone"
$ws2.Range("E1").Value = "aadhaar"
$ws2.Range("F1").Value = "pan"

$ws2.Range("A2").Value = "Shekhar"
$ws2.Range("B2").Value = "Palo"
$ws2.Range("C2").Value = "shekhar123@test.com"
$ws2.Range("D2").Value = 8896147850
$ws2.Range("E2").Value = 987286540931
$ws2.Range("F2").Value = "CHUPK3570Q"

$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:shekhar123@test.com")

# column widths for GTM sheet
# NOTE: the engine quantizes ColumnWidth to 1/6-character steps, so these
# inputs are chosen to land on the closest achievable stored width to the
# target (11.28515625, 12.7109375, 25.28515625, 15.5703125, 25.85546875,
# 20.28515625 "char width units" respectively).
$ws2.Columns.Item(1).ColumnWidth = 10.5
$ws2.Columns.Item(2).ColumnWidth = 11.83
$ws2.Columns.Item(3).ColumnWidth = 24.5
$ws2.Columns.Item(4).ColumnWidth = 14.665
$ws2.Columns.Item(5).ColumnWidth = 25.0
$ws2.Columns.Item(6).ColumnWidth = 19.5

$ws2.Activate()
$ws2.Select()

$ws1.Range("B2").Select()
$ws2.Range("G2").Select()
